$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.346
$ws.Range("B3").Value = 6.798999999999999
$ws.Range("B5").Value = 6.243
$ws.Range("E7").Value = 13.078
$ws.Range("A9").Value = -21.096
$ws.Range("E9").Value = 12.889
$ws.Range("B11").Value = 6.93
$ws.Range("B12").Value = 6.582000000000001
$ws.Range("A13").Value = -21.832
$ws.Range("A16").Value = -20.763
$ws.Range("A18").Value = -21.868
$ws.Range("A20").Value = -21.605
$ws.Range("B21").Value = 6.928999999999999
$ws.Range("E21").Value = 12.608
